# Apply updates to data/correlation_analysis.xlsx as described in the
# "updated stuff from replication-package branch" commit.
#
# Sheet "all_tools" (first sheet): column J width changes, and several
# correlation-statistics values in rows 9-12 are updated.
#
# Sheet "openjml" (fifth sheet): row 9 values updated (F,G,I,J,K,L) and
# rows 10-12 get updated F,G values only.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: all_tools
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("all_tools")

# Column J (10th column) width: 21.7109375 -> 20.7109375
# Note: this runtime's ColumnWidth setter snaps the stored OOXML width to
# the nearest 1/6-character (whole-pixel) grid point, so 20.7109375 itself
# is not landable exactly; 19.8333... is the input that lands on the
# closest reachable grid value (20.6666...) to the target.
$wsAll.Columns.Item(10).ColumnWidth = 19.8333333333333

# Row 9
$wsAll.Range("F9").Value = 94
$wsAll.Range("G9").Value = 1200
$wsAll.Range("I9").Value = -0.1679930849817192
$wsAll.Range("J9").Value = 0.01571209940531982
$wsAll.Range("K9").Value = -0.2581562594294852
$wsAll.Range("L9").Value = 0.009508181415596517

# Row 10
$wsAll.Range("G10").Value = 863
$wsAll.Range("I10").Value = -0.01788960397609135
$wsAll.Range("J10").Value = 0.8831700141519032
$wsAll.Range("K10").Value = -0.01605248793883186
$wsAll.Range("L10").Value = 0.9118983828429555

# Row 11
$wsAll.Range("G11").Value = 863
$wsAll.Range("I11").Value = -0.0319066731228802
$wsAll.Range("J11").Value = 0.7837672300250216
$wsAll.Range("K11").Value = -0.03045075791865187
$wsAll.Range("L11").Value = 0.8337285990222609

# Row 12
$wsAll.Range("G12").Value = 863
$wsAll.Range("I12").Value = -0.2441570787714312
$wsAll.Range("J12").Value = 0.03455416553919029
$wsAll.Range("K12").Value = -0.3094340106343606
$wsAll.Range("L12").Value = 0.02876701130235894

# ---------------------------------------------------------------------
# Sheet: openjml
# ---------------------------------------------------------------------
$wsJml = $wb.Worksheets.Item("openjml")

# Row 9
$wsJml.Range("F9").Value = 69
$wsJml.Range("G9").Value = 808
$wsJml.Range("I9").Value = -0.1508932577851104
$wsJml.Range("J9").Value = 0.03976539977998658
$wsJml.Range("K9").Value = -0.2101007063448729
$wsJml.Range("L9").Value = 0.03589890255176344

# Row 10
$wsJml.Range("F10").Value = 41
$wsJml.Range("G10").Value = 219

# Row 11
$wsJml.Range("F11").Value = 41
$wsJml.Range("G11").Value = 219

# Row 12
$wsJml.Range("F12").Value = 41
$wsJml.Range("G12").Value = 219
